# Apply updated crypto price/volume values (and a few coin name/link changes
# for rows 49-51) as described in the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'29.311.85"
$ws.Range("E2").Value = "'  +1.47%  "
$ws.Range("D3").Value = "'1.841.40"
$ws.Range("E3").Value = "'  +0.55%  "
$ws.Range("D4").Value = "'0.9995"
$ws.Range("E4").Value = "'  +0.13%  "
$ws.Range("D5").Value = "'242.89"
$ws.Range("E5").Value = "'  -0.67%  "
$ws.Range("D6").Value = "'0.6866"
$ws.Range("E6").Value = "'  -0.69%  "
$ws.Range("D7").Value = "'1.000"
$ws.Range("E7").Value = "'  +0.12%  "
$ws.Range("D8").Value = "'0.3029"
$ws.Range("E8").Value = "'  -0.51%  "
$ws.Range("D9").Value = "'0.07530"
$ws.Range("E9").Value = "'  -1.84%  "
$ws.Range("D10").Value = "'23.27"
$ws.Range("E10").Value = "'  +0.23%  "
$ws.Range("D11").Value = "'0.07659"
$ws.Range("E11").Value = "'  -1.58%  "
$ws.Range("D12").Value = "'1.836.73"
$ws.Range("E12").Value = "'  +0.27%  "
$ws.Range("D13").Value = "'5.085"
$ws.Range("E13").Value = "'  -0.02%  "
$ws.Range("D14").Value = "'0.6867"
$ws.Range("E14").Value = "'  +0.94%  "
$ws.Range("D15").Value = "'88.94"
$ws.Range("E15").Value = "'  -4.33%  "
$ws.Range("D16").Value = "'6.288"
$ws.Range("E16").Value = "'  -2.71%  "
$ws.Range("D17").Value = "'29.312.78"
$ws.Range("E17").Value = "'  +1.54%  "
$ws.Range("D18").Value = "'0.000008237"
$ws.Range("E18").Value = "'  -0.25%  "
$ws.Range("D19").Value = "'2.085.36"
$ws.Range("E19").Value = "'  +0.74%  "
$ws.Range("D20").Value = "'232.81"
$ws.Range("E20").Value = "'  -3.77%  "
$ws.Range("D21").Value = "'12.61"
$ws.Range("E21").Value = "'  -0.42%  "
$ws.Range("D22").Value = "'0.9999"
$ws.Range("E22").Value = "'  +0.11%  "
$ws.Range("E23").Value = "'  +0.68%  "
$ws.Range("D24").Value = "'1.0000"
$ws.Range("E24").Value = "'  +0.06%  "
$ws.Range("D25").Value = "'0.1462"
$ws.Range("E25").Value = "'  -1.98%  "
$ws.Range("D26").Value = "'160.20"
$ws.Range("E26").Value = "'  -0.24%  "
$ws.Range("D27").Value = "'8.838"
$ws.Range("E27").Value = "'  +1.29%  "
$ws.Range("D28").Value = "'18.08"
$ws.Range("E28").Value = "'  -0.63%  "
$ws.Range("D29").Value = "'1.519"
$ws.Range("E29").Value = "'  -1.30%  "
$ws.Range("D30").Value = "'4.231"
$ws.Range("E30").Value = "'  +0.16%  "
$ws.Range("E31").Value = "'  -0.23%  "
$ws.Range("E32").Value = "'  +1.22%  "
$ws.Range("D33").Value = "'0.05155"
$ws.Range("E33").Value = "'  +0.92%  "
$ws.Range("D34").Value = "'0.7748"
$ws.Range("E34").Value = "'  +0.15%  "
$ws.Range("D35").Value = "'1.845"
$ws.Range("E35").Value = "'  -0.54%  "
$ws.Range("D36").Value = "'1.140"
$ws.Range("E36").Value = "'  +0.06%  "
$ws.Range("D37").Value = "'2.673"
$ws.Range("E37").Value = "'  -0.61%  "
$ws.Range("D38").Value = "'1.289.39"
$ws.Range("E38").Value = "'  +1.78%  "
$ws.Range("D39").Value = "'0.01845"
$ws.Range("E39").Value = "'  -0.70%  "
$ws.Range("D40").Value = "'2.701"
$ws.Range("E40").Value = "'  +0.06%  "
$ws.Range("D41").Value = "'0.9448"
$ws.Range("E41").Value = "'  -1.70%  "
$ws.Range("D42").Value = "'105.70"
$ws.Range("E42").Value = "'  -0.98%  "
$ws.Range("D43").Value = "'0.9996"
$ws.Range("E43").Value = "'  +0.03%  "
$ws.Range("D44").Value = "'5.693"
$ws.Range("E44").Value = "'  -6.44%  "
$ws.Range("D45").Value = "'9.675"
$ws.Range("E45").Value = "'  +0.31%  "
$ws.Range("D46").Value = "'1.985.65"
$ws.Range("E46").Value = "'  +0.65%  "
$ws.Range("D47").Value = "'0.5204"
$ws.Range("E47").Value = "'  +0.83%  "
$ws.Range("E48").Value = "'  +1.22%  "
$ws.Range("B49").Value = "'Aave"
$ws.Range("C49").Value = "'https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D49").Value = "'63.36"
$ws.Range("E49").Value = "'  -0.71%  "
$ws.Range("B50").Value = "'Cronos"
$ws.Range("C50").Value = "'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D50").Value = "'0.05925"
$ws.Range("E50").Value = "'  +0.69%  "
$ws.Range("B51").Value = "'Aptos"
$ws.Range("C51").Value = "'https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D51").Value = "'6.905"
$ws.Range("E51").Value = "'  -0.55%  "

# Reset number format/style on touched cells so forcing text values via the
# leading apostrophe does not leave a stray quote-prefix style behind.
$ws.Range("B2:E51").Style = "Normal"
